# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap country names / stats (rows 192 & 193): Belice <-> Nueva Caledonia ---
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("A193").Value = "Belice"

$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# --- Swap country names / stats (rows 212 & 213): Butan <-> Islas Virgenes Britanicas ---
$ws.Range("A212").Value = "Islas Virgenes Britanicas"
$ws.Range("A213").Value = "Butan"

$ws.Range("D212").Value = 4
$ws.Range("H212").Value = 1
$ws.Range("D213").Value = 5
$ws.Range("H213").Value = 0

# --- Estados Unidos (row 4) updated totals ---
$ws.Range("B4").Value = 1365308
$ws.Range("C4").Value = 17999
$ws.Range("D4").Value = 240853
$ws.Range("E4").Value = 1043738
$ws.Range("G4").Value = 680
$ws.Range("H4").Value = 80717

# --- Brasil (row 11) updated totals ---
$ws.Range("B11").Value = 162699
$ws.Range("C11").Value = 6638
$ws.Range("E11").Value = 89891
$ws.Range("G11").Value = 467
$ws.Range("H11").Value = 11123

# --- Arabia Saudita (row 20) updated totals ---
$ws.Range("F20").Value = 143
